$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  41"
$ws.Range("C9").Value = "Report Covering the Week  10/6/2025  Through  10/12/2025"

# --- C22: number -> text "0" (shared string), preserving style 13 ---
$fmtSrc = $ws.Range("D23")
$c22 = $ws.Range("C22")
$fmtSrc.Copy()
$c22.PasteSpecial(-4122)   # xlPasteFormats: copies style (incl. number format) from D23 (s=13)
$fmtSrc.Copy()
$c22.PasteSpecial(-4163)   # xlPasteValues: copies the text value "0" (shared string) from D23
$excel.CutCopyMode = $false

# --- Numeric cell updates ---
# Row 14
$ws.Range("G14").Value = 1
$ws.Range("M14").Value = -40

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("I15").Value = 30
$ws.Range("J15").Value = 28
$ws.Range("K15").Value = 7.142857142857
$ws.Range("L15").Value = 87.5
$ws.Range("M15").Value = 42.857142857142
$ws.Range("N15").Value = -55.223880597014

# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 12.5
$ws.Range("F16").Value = 36
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 330
$ws.Range("J16").Value = 309
$ws.Range("K16").Value = 6.796116504854
$ws.Range("L16").Value = 22.222222222222
$ws.Range("M16").Value = -3.790087463556
$ws.Range("N16").Value = -80.450236966824

# Row 17
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -23.076923076923
$ws.Range("F17").Value = 65
$ws.Range("G17").Value = 69
$ws.Range("H17").Value = -5.797101449275
$ws.Range("I17").Value = 677
$ws.Range("J17").Value = 700
$ws.Range("K17").Value = -3.285714285714
$ws.Range("L17").Value = 5.124223602484
$ws.Range("M17").Value = 99.117647058823
$ws.Range("N17").Value = -21.643518518518

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 13
$ws.Range("H18").Value = -27.777777777777
$ws.Range("I18").Value = 148
$ws.Range("J18").Value = 170
$ws.Range("K18").Value = -12.941176470588
$ws.Range("L18").Value = 10.447761194029
$ws.Range("M18").Value = -13.45029239766
$ws.Range("N18").Value = -91.764051196438

# Row 19
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 60
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 25
$ws.Range("I19").Value = 523
$ws.Range("J19").Value = 537
$ws.Range("K19").Value = -2.607076350093
$ws.Range("L19").Value = 19.134396355353
$ws.Range("M19").Value = 118.828451882845
$ws.Range("N19").Value = -33.460559796437

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -44.444444444444
$ws.Range("F20").Value = 18
$ws.Range("H20").Value = -10
$ws.Range("I20").Value = 200
$ws.Range("J20").Value = 173
$ws.Range("K20").Value = 15.606936416185
$ws.Range("L20").Value = -11.504424778761
$ws.Range("M20").Value = 106.185567010309
$ws.Range("N20").Value = -71.791255289139

# Row 21
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -2.272727272727
$ws.Range("F21").Value = 193
$ws.Range("G21").Value = 190
$ws.Range("H21").Value = 1.578947368421
$ws.Range("I21").Value = 1914
$ws.Range("J21").Value = 1941
$ws.Range("K21").Value = -1.391035548686
$ws.Range("L21").Value = 9.936817920735
$ws.Range("M21").Value = 56.756756756756
$ws.Range("N21").Value = -67.885906040268

# Row 22
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = 26.315789473684
$ws.Range("L22").Value = 4.347826086956
$ws.Range("M22").Value = -7.692307692307

# Row 23
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 31
$ws.Range("K23").Value = 29.166666666666
$ws.Range("L23").Value = 34.782608695652
$ws.Range("M23").Value = 63.157894736842

# Row 24
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 12.5
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = 9.859154929577
$ws.Range("I24").Value = 795
$ws.Range("J24").Value = 751
$ws.Range("K24").Value = 5.858854860186
$ws.Range("L24").Value = -3.985507246376
$ws.Range("M24").Value = 8.310626702997

# Row 25
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -85.714285714285
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = -77.777777777777
$ws.Range("I25").Value = 212
$ws.Range("J25").Value = 266
$ws.Range("K25").Value = -20.300751879699
$ws.Range("L25").Value = 23.255813953488

# Row 26
$ws.Range("C26").Value = 23
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 43.75
$ws.Range("F26").Value = 84
$ws.Range("G26").Value = 79
$ws.Range("H26").Value = 6.32911392405
$ws.Range("I26").Value = 826
$ws.Range("J26").Value = 898
$ws.Range("K26").Value = -8.017817371937
$ws.Range("L26").Value = -16.649848637739
$ws.Range("M26").Value = -19.961240310077

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -87.5
$ws.Range("I27").Value = 42
$ws.Range("J27").Value = 41
$ws.Range("K27").Value = 2.439024390243
$ws.Range("L27").Value = -8.695652173913

# Row 28
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 66.666666666666
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 75
$ws.Range("J28").Value = 77
$ws.Range("K28").Value = -2.597402597402
$ws.Range("L28").Value = 7.142857142857

# Row 29
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = -85.714285714285
$ws.Range("J29").Value = 56
$ws.Range("K29").Value = -58.928571428571
$ws.Range("L29").Value = -54.901960784313
$ws.Range("M29").Value = -53.061224489795
$ws.Range("N29").Value = -84.027777777777

# Row 30
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -75
$ws.Range("J30").Value = 43
$ws.Range("K30").Value = -55.813953488372
$ws.Range("L30").Value = -52.5
$ws.Range("M30").Value = -54.761904761904
$ws.Range("N30").Value = -85.60606060606
